# Updates the "alloptions" sheet: rewrites the method-locator header/
# description blocks and adds two new imported-library columns (Wait, Title).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: column headers (existing columns kept, two new ones appended) ---
$ws.Range("A1").Value = "Browser"
$ws.Range("B1").Value = "URL"
$ws.Range("C1").Value = "Search"
$ws.Range("D1").Value = "Button"
$ws.Range("E1").Value = "End"
$ws.Range("F1").Value = "Back"
$ws.Range("G1").Value = "Forward"
$ws.Range("H1").Value = "Clear"
$ws.Range("I1").Value = "Wait"
$ws.Range("J1").Value = "Title"

# --- Row 2: per-column input hints ---
$ws.Range("A2").Value = "Browser input"
$ws.Range("B2").Value = "Specific URL"
$ws.Range("C2").Value = "Specific locator"
$ws.Range("D2").Value = "Specific locator"
$ws.Range("F2").Value = "Number of"
$ws.Range("G2").Value = "Number of"
$ws.Range("I2").Value = "Seconds"
$ws.Range("J2").Value = "String"

# --- Row 3 (unchanged) ---
$ws.Range("C3").Value = "Location"
$ws.Range("D3").Value = "Location"

# --- Row 4: keep C4, clear the old H4 "Condition:" (moved to row 8) ---
$ws.Range("C4").Value = "Search input"
$ws.Range("H4:H5").ClearContents()

# --- Row 6: "Description:" banner across all 10 columns ---
$ws.Range("A6:J6").Value = "Description:"

# --- Row 7: per-column descriptions for the new/notable columns ---
$ws.Range("A7").Value = "A1:Specification of step A2:Which browser shoud be used"
$ws.Range("B7").Value = "B1: Specification of step B2: Input of URL"
$ws.Range("H7").Value = "Clear string in search box"
$ws.Range("I7").Value = "Waiting  for some time(in seconds)"
$ws.Range("J7").Value = "Check if title matches with input string"

# --- Row 8: "Inputs:"/"Condition:" banner ---
$ws.Range("A8").Value = "Inputs:"
$ws.Range("B8:J8").Value = "Condition:"

# --- Row 9: supplemental notes ---
$ws.Range("A9").Value = "(Chrome, Firefox, Safari, Edge)"
$ws.Range("B9").Value = "Anything on web"
$ws.Range("H9").Value = "Must be after search"
$ws.Range("J9").Value = "-"

# --- View state: scroll right a bit and move the active selection ---
$ws.Range("K1").Select()
